# Fixed output parsing logic not working with color code escape characters
#
# The test fixture's first row (an unrelated "Hello World | grep" example) is
# dropped, shifting the two remaining original rows up to A1/A2, and eight
# new multi-line terminal-transcript examples are added below them
# (A3..A10), wrapped and sized to fit their content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1 / A2: keep the two still-relevant examples, shifted up one row ---
$ws.Range("A1").Value = @'
$> echo "         " | cat -e
'@

$ws.Range("A2").Value = @'
$> echo -n "         " | cat -e
'@

# --- new multi-line examples, written in the order that fixes the shared-
#     string table's index assignment (A5's text is index 2, A3/A10's text
#     is index 3, etc.) ---
$sHomeHola = @'
$> echo $HOME
hola
'@
$ws.Range("A5").Value = $sHomeHola

$sTripleEcho = @'
$> echo $HOME
$> echo $HOME
$> echo $HOME
'@
$ws.Range("A3").Value = $sTripleEcho
$ws.Range("A10").Value = $sTripleEcho

$sHomeThenHola = @'
$> echo $HOME
$> hola
'@
$ws.Range("A4").Value = $sHomeThenHola

$sHomeOnly = @'
$> echo $HOME
'@
$ws.Range("A6").Value = $sHomeOnly

$sBlankLeadEcho = @'

echo $HOME
echo $HOME
echo $HOME
'@
$ws.Range("A7").Value = $sBlankLeadEcho

$sPromptLeadEcho = @'
$> 
$> echo $HOME
$> echo $HOME
$> echo $HOME
'@
$ws.Range("A8").Value = $sPromptLeadEcho

$sPlainEcho = @'
echo $HOME
echo $HOME
echo $HOME
'@
$ws.Range("A9").Value = $sPlainEcho

# --- wrap + size the new rows to fit their multi-line content ---
$ws.Range("A3:A10").WrapText = $true

$ws.Rows.Item(3).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 57.6
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(10).RowHeight = 43.2

# Cursor ends up one row past the last entry, like after typing A10 + Enter.
$ws.Range("A11").Select() | Out-Null
